$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027438755998537
$ws.Range("D2").Value = 1.032276575641003
$ws.Range("E2").Value = 1.027518295290458
$ws.Range("F2").Value = 1.03795793616981
$ws.Range("I2").Value = 1.034642685647126
$ws.Range("J2").Value = 1.032596626213845
$ws.Range("K2").Value = 1.03508227233661
$ws.Range("L2").Value = 1.030337790339339
$ws.Range("M2").Value = 1.040747337688945
$ws.Range("N2").Value = 1.014821978095935
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028535683675962
$ws.Range("D3").Value = 1.032806439804623
$ws.Range("E3").Value = 1.028454351913477
$ws.Range("F3").Value = 1.039245160370384
$ws.Range("I3").Value = 1.034871073262916
$ws.Range("J3").Value = 1.033332981854826
$ws.Range("K3").Value = 1.035421621698634
$ws.Range("L3").Value = 1.031081242145633
$ws.Range("M3").Value = 1.041843214701677
$ws.Range("N3").Value = 1.015071526364898
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029245378401928
$ws.Range("D4").Value = 1.033149408415588
$ws.Range("E4").Value = 1.029060308299605
$ws.Range("F4").Value = 1.040078280243672
$ws.Range("I4").Value = 1.035017657864211
$ws.Range("J4").Value = 1.033808832603767
$ws.Range("K4").Value = 1.035640603406677
$ws.Range("L4").Value = 1.031561950383285
$ws.Range("M4").Value = 1.042551974951399
$ws.Range("N4").Value = 1.0152326484628
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.029543713293554
$ws.Range("D5").Value = 1.033293617897597
$ws.Range("E5").Value = 1.02931511598466
$ws.Range("F5").Value = 1.040428573349856
$ws.Range("I5").Value = 1.035078995260182
$ws.Range("J5").Value = 1.034008732060632
$ws.Range("K5").Value = 1.035732519122288
$ws.Range("L5").Value = 1.031763955428717
$ws.Range("M5").Value = 1.042849856424026
$ws.Range("N5").Value = 1.015300299862757
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.029593803844919
$ws.Range("D6").Value = 1.033317832765773
$ws.Range("E6").Value = 1.029357903050991
$ws.Range("F6").Value = 1.040487392056394
$ws.Range("I6").Value = 1.0350892772597
$ws.Range("J6").Value = 1.034042287379461
$ws.Range("K6").Value = 1.035747943703756
$ws.Range("L6").Value = 1.031797868004985
$ws.Range("M6").Value = 1.042899867304188
$ws.Range("N6").Value = 1.015311653885784
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029249364847742
$ws.Range("D7").Value = 1.033151335251045
$ws.Range("E7").Value = 1.029063712802205
$ws.Range("F7").Value = 1.040082960680633
$ws.Range("I7").Value = 1.03501847858443
$ws.Range("J7").Value = 1.033811504251946
$ws.Range("K7").Value = 1.035641832156248
$ws.Range("L7").Value = 1.031564649917095
$ws.Range("M7").Value = 1.042555955576103
$ws.Range("N7").Value = 1.015233552755458
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027809487236024
$ws.Range("D8").Value = 1.032455622052681
$ws.Range("E8").Value = 1.027834585191538
$ws.Range("F8").Value = 1.038392919275826
$ws.Range("I8").Value = 1.034720118194936
$ws.Range("J8").Value = 1.032845609772075
$ws.Range("K8").Value = 1.035197080952014
$ws.Range("L8").Value = 1.030589116966152
$ws.Range("M8").Value = 1.041117766807599
$ws.Range("N8").Value = 1.014906387061784
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.025271499318657
$ws.Range("D9").Value = 1.031230585870206
$ws.Range("E9").Value = 1.025670729188179
$ws.Range("F9").Value = 1.035416304174149
$ws.Range("I9").Value = 1.034185201992948
$ws.Range("J9").Value = 1.031138810194623
$ws.Range("K9").Value = 1.034408800991175
$ws.Range("L9").Value = 1.028867374695085
$ws.Range("M9").Value = 1.038580787219087
$ws.Range("N9").Value = 1.014327180014466
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.02357894609643
$ws.Range("D10").Value = 1.030414564591519
$ws.Range("E10").Value = 1.024229506680822
$ws.Range("F10").Value = 1.033432762930697
$ws.Range("I10").Value = 1.03382243202435
$ws.Range("J10").Value = 1.029997710371601
$ws.Range("K10").Value = 1.033880238689525
$ws.Range("L10").Value = 1.027717693789551
$ws.Range("M10").Value = 1.036887559086627
$ws.Range("N10").Value = 1.013939224769683
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022845903678404
$ws.Range("D11").Value = 1.030061390572328
$ws.Range("E11").Value = 1.023605757475677
$ws.Range("F11").Value = 1.032574048828574
$ws.Range("I11").Value = 1.033663888254134
$ws.Range("J11").Value = 1.029502828023777
$ws.Range("K11").Value = 1.033650650065717
$ws.Range("L11").Value = 1.027219425010816
$ws.Range("M11").Value = 1.036153899594326
$ws.Range("N11").Value = 1.013770803996624
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022573594293452
$ws.Range("D12").Value = 1.029930232398339
$ws.Range("E12").Value = 1.023374115349268
$ws.Range("F12").Value = 1.032255107826408
$ws.Range("I12").Value = 1.033604778373731
$ws.Range("J12").Value = 1.029318889048072
$ws.Range("K12").Value = 1.033565263273672
$ws.Range("L12").Value = 1.027034277783571
$ws.Range("M12").Value = 1.035881311744173
$ws.Range("N12").Value = 1.013708179767233
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022632006725287
$ws.Range("D13").Value = 1.029958365061205
$ws.Range("E13").Value = 1.023423801292347
$ws.Range("F13").Value = 1.032323520719898
$ws.Range("I13").Value = 1.033617467584348
$ws.Range("J13").Value = 1.029358349916739
$ws.Range("K13").Value = 1.033583583881015
$ws.Range("L13").Value = 1.027073995578401
$ws.Range("M13").Value = 1.035939786143984
$ws.Range("N13").Value = 1.013721615831795
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022823395013547
$ws.Range("D14").Value = 1.0300505484525
$ws.Range("E14").Value = 1.02358660892725
$ws.Range("F14").Value = 1.032547684596378
$ws.Range("I14").Value = 1.033659006694192
$ws.Range("J14").Value = 1.029487625974263
$ws.Range("K14").Value = 1.033643594154177
$ws.Range("L14").Value = 1.027204122073726
$ws.Range("M14").Value = 1.036131368918891
$ws.Range("N14").Value = 1.013765628790136
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022941312327472
$ws.Range("D15").Value = 1.030107349209566
$ws.Range("E15").Value = 1.023686926180606
$ws.Range("F15").Value = 1.03268580242568
$ws.Range("I15").Value = 1.033684571202409
$ws.Range("J15").Value = 1.029567261612055
$ws.Range("K15").Value = 1.03368055425237
$ws.Range("L15").Value = 1.027284288274062
$ws.Range("M15").Value = 1.036249399521927
$ws.Range("N15").Value = 1.013792737976006
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023627592307575
$ws.Range("D16").Value = 1.030438007230005
$ws.Range("E16").Value = 1.024270909349456
$ws.Range("F16").Value = 1.03348975639724
$ws.Range("I16").Value = 1.033832923248384
$ws.Range("J16").Value = 1.03003053762262
$ws.Range("K16").Value = 1.033895460649924
$ws.Range("L16").Value = 1.027750752736747
$ws.Range("M16").Value = 1.036936239368373
$ws.Range("N16").Value = 1.013950393156813
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024058034934453
$ws.Range("D17").Value = 1.030645466150968
$ws.Range("E17").Value = 1.024637309169413
$ws.Range("F17").Value = 1.033994100042937
$ws.Range("I17").Value = 1.033925589140774
$ws.Range("J17").Value = 1.030320929506277
$ws.Range("K17").Value = 1.034030073947544
$ws.Range("L17").Value = 1.028043232907768
$ws.Range("M17").Value = 1.037366945897422
$ws.Range("N17").Value = 1.014049169922451
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024309089907024
$ws.Range("D18").Value = 1.030766489624822
$ws.Range("E18").Value = 1.024851053890691
$ws.Range("F18").Value = 1.034288291938128
$ws.Range("I18").Value = 1.033979498529485
$ws.Range("J18").Value = 1.030490235100457
$ws.Range("K18").Value = 1.034108522331875
$ws.Range("L18").Value = 1.028213788320552
$ws.Range("M18").Value = 1.037618123488371
$ws.Range("N18").Value = 1.014106742911021
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024394690675125
$ws.Range("D19").Value = 1.030807758216657
$ws.Range("E19").Value = 1.024923940402183
$ws.Range("F19").Value = 1.03438860668608
$ws.Range("I19").Value = 1.033997856312483
$ws.Range("J19").Value = 1.030547951218308
$ws.Range("K19").Value = 1.034135259439231
$ws.Range("L19").Value = 1.028271935975853
$ws.Range("M19").Value = 1.037703760769031
$ws.Range("N19").Value = 1.014126366721021
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02401185406834
$ws.Range("D20").Value = 1.030623206089314
$ws.Range("E20").Value = 1.02459799484595
$ws.Range("F20").Value = 1.033939987023142
$ws.Range("I20").Value = 1.033915661557933
$ws.Range("J20").Value = 1.030289780975756
$ws.Range("K20").Value = 1.034015638376676
$ws.Range("L20").Value = 1.028011857037998
$ws.Range("M20").Value = 1.037320739959015
$ws.Range("N20").Value = 1.014038576437509
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022767036628513
$ws.Range("D21").Value = 1.030023402004677
$ws.Range("E21").Value = 1.023538664867433
$ws.Range("F21").Value = 1.032481673287076
$ws.Range("I21").Value = 1.033646780525663
$ws.Range("J21").Value = 1.029449560650131
$ws.Range("K21").Value = 1.033625925584245
$ws.Range("L21").Value = 1.027165804945861
$ws.Range("M21").Value = 1.03607495462575
$ws.Range("N21").Value = 1.013752669873181
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021984226363776
$ws.Range("D22").Value = 1.029646434523774
$ws.Range("E22").Value = 1.022872888579524
$ws.Range("F22").Value = 1.031564909277869
$ws.Range("I22").Value = 1.033476453425046
$ws.Range("J22").Value = 1.028920599841848
$ws.Range("K22").Value = 1.033380276790101
$ws.Range("L22").Value = 1.026633464635589
$ws.Range("M22").Value = 1.035291250600779
$ws.Range("N22").Value = 1.013572531464187
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022399223026378
$ws.Range("D23").Value = 1.029846257250671
$ws.Range("E23").Value = 1.023225804065692
$ws.Range("F23").Value = 1.032050891110472
$ws.Range("I23").Value = 1.033566867519242
$ws.Range("J23").Value = 1.029201076741743
$ws.Range("K23").Value = 1.033510558562943
$ws.Range("L23").Value = 1.026915705801908
$ws.Range("M23").Value = 1.03570674824871
$ws.Range("N23").Value = 1.01366806207386
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024032721259619
$ws.Range("D24").Value = 1.030633264403055
$ws.Range("E24").Value = 1.024615759203843
$ws.Range("F24").Value = 1.033964438315073
$ws.Range("I24").Value = 1.033920147841074
$ws.Range("J24").Value = 1.030303855887657
$ws.Range("K24").Value = 1.034022161403837
$ws.Range("L24").Value = 1.028026034576632
$ws.Range("M24").Value = 1.037341618576026
$ws.Range("N24").Value = 1.014043363306536
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025927724264322
$ws.Range("D25").Value = 1.031547173799682
$ws.Range("E25").Value = 1.026229899273417
$ws.Range("F25").Value = 1.036185669425935
$ws.Range("I25").Value = 1.034324576556522
$ws.Range("J25").Value = 1.031580626776469
$ws.Range("K25").Value = 1.03461312876387
$ws.Range("L25").Value = 1.029312811553763
$ws.Range("M25").Value = 1.03923698715324
$ws.Range("N25").Value = 1.014477239036516

Write-Output "Updated 264 cells"